# Append the new match row (row 74) to the HNL 2023-2024 sheet,
# mirroring the existing data rows (values + per-column formatting).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 74
$prev = $row - 1

$ws.Cells.Item($row, 1).Value = 73
$ws.Cells.Item($row, 2).Value = "croatia"
$ws.Cells.Item($row, 3).Value = "hnl"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45254.75
$ws.Cells.Item($row, 6).Value = "Lok. Zagreb"
$ws.Cells.Item($row, 7).Value = 1
$ws.Cells.Item($row, 8).Value = "Slaven Belupo"
$ws.Cells.Item($row, 9).Value = 3
$ws.Cells.Item($row, 10).Value = 1.73
$ws.Cells.Item($row, 11).Value = "19/11/2023 10:12"
$ws.Cells.Item($row, 12).Value = 1.83
$ws.Cells.Item($row, 13).Value = "24/11/2023 17:55"
$ws.Cells.Item($row, 14).Value = 3.67
$ws.Cells.Item($row, 15).Value = "19/11/2023 10:12"
$ws.Cells.Item($row, 16).Value = 3.59
$ws.Cells.Item($row, 17).Value = "24/11/2023 17:58"
$ws.Cells.Item($row, 18).Value = 4.81
$ws.Cells.Item($row, 19).Value = "19/11/2023 10:12"
$ws.Cells.Item($row, 20).Value = 4.47
$ws.Cells.Item($row, 21).Value = "24/11/2023 17:55"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/croatia/hnl/lok-zagreb-slaven-belupo/8f55XWPA/"

# Copy the formatting (styles) of the previous data row onto the new row,
# cell by cell, so the new row reuses the same style indexes (bold/border
# index column, custom date format column, plain cells elsewhere) instead
# of minting new styles.
for ($col = 1; $col -le 22; $col++) {
    $ws.Cells.Item($prev, $col).Copy()
    $ws.Cells.Item($row, $col).PasteSpecial(-4122)
}

$excel.CutCopyMode = $false
